# position sync create a building by player's position
#
# The skill cooldown (column P, "CoolDownTime") drops from 2 to 1 for the
# rows whose cooldown was 2, and the damage distance (column R,
# "DamageDistance") drops from 2.5 to 2 for every data row (11-61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose CoolDownTime (column P) goes from 2 -> 1.
$coolDownRows = @(11, 12, 13, 14, 15, 16, 17, 18, 19, 38, 41, 44, 47, 50, 53, 56, 59)
foreach ($r in $coolDownRows) {
    $ws.Cells.Item($r, 16).Value = 1
}

# DamageDistance (column R) goes from 2.5 -> 2 for every data row (11-61).
for ($r = 11; $r -le 61; $r++) {
    $ws.Cells.Item($r, 18).Value = 2
}
